# am modificat la arhitectura
# Applies the Architect. Design Phase Defects sheet edits:
#  - tightened/fixed a couple of review comments (A01, A02, A06, A07 rows)
#  - added 4 new rows of comments (A08-A11-ish, rows 21-24) describing the
#    architecture fixes that were made
#  - cleared the leftover border/italic styling on E20 now that it is the
#    last "empty" comment row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Architect. Design Phase Defects")

# --- Row 10 (A01 comment): drop the now-resolved sub-questions ---
$ws.Range("E10").Value = "BibliotecaCtrl: ce cauta acol modifica carte  si sterge carte? CartiRepo: ce cauta acol toate functiile ce trebuiau sa fie in controller si unde is functiile de baza a unui repo? Validator nu o fost specificat in requirements. Carte: constructorul nu are parametrii. Operatii cu referenti si cuvinte cheie care nush ce cauta acolo. Cautarile iara nu au ce cauta acolo. Parsarea din string n-o fost specificata nicaria. "

# --- Row 11 (A02 comment): word choice fix ---
$ws.Range("E11").Value = "no, separate things better and don’t add unneeded stuff(see above)"

# --- Row 15 (A06 comment): typo fix ---
$ws.Range("E15").Value = "well, there is inversion of control, but that’s it"

# --- Row 16 (A07 comment): clarified wording ---
$ws.Range("E16").Value = "“Consola” is not clear, would be better named “UIConsola” or something"

# --- Rows 21-24: new follow-up comments about the architecture fixes ---
$ws.Range("E21").Value = "adaugat editura la carte"
$ws.Range("E22").Value = "scos modifica si sterge carte din controller"
$ws.Range("E23").Value = "scos functiile irelevante din repository si le-am pus in controller"
$ws.Range("E24").Value = "scos functiile de cautare dupa autor si cuvinte cheie din Carte"

# --- Row 20 reverts to the plain/unstyled look (no border, no italic) ---
$ws.Range("E20").Borders.LineStyle = 0
$ws.Range("E20").Font.Italic = $false

# --- Row heights: shrink row 10 (shorter text) and tidy 20-24 ---
$ws.Rows.Item(10).RowHeight = 102.2
$ws.Rows.Item(20).RowHeight = 13.8
$ws.Rows.Item(21).RowHeight = 13.8
$ws.Rows.Item(22).RowHeight = 13.8
$ws.Rows.Item(23).RowHeight = 23.85
$ws.Rows.Item(24).RowHeight = 23.85

# --- Cursor/selection position as left by the author ---
$ws.Application.ActiveWindow.ScrollRow = 9
$ws.Range("E24").Select()

$ws3 = $wb.Worksheets.Item("Coding Phase Defects")
$ws3.Range("E14").Select()
